# Applies the changes described by the upstream diff:
#  - G2 formula: (D2/C2*100)/100  ->  ((K2/114)*100)/100
#  - G3:G21 (shared formula) : (D3/C3*100)/100 -> ((K3/114)*100)/100
#  - Selection / scroll position moved from O1 to P8 (window scrolled right to show column H onward)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "Percentual de aproveitamento" column formulas (column G) ---
# G2 keeps its own (non-shared) formula, matching the original layout.
$ws.Range("G2").Formula = "=((K2/114)*100)/100"

# G3:G21 is filled from a single formula so Excel keeps it as one shared-formula group,
# exactly like the original (D3/C3*100)/100 group that spanned G3:G21.
$ws.Range("G3:G21").Formula = "=((K3/114)*100)/100"

# --- Update the view: scroll the sheet so column H is the first visible column, ---
# --- and move the active selection to P8 (previously O1). ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 8
$ws.Range("P8").Select()
